$d = $word.ActiveDocument

# The CV lists several papers whose citation ends with
# ", draft available on request". We only want to drop that phrase from
# the entry for "Coded many-user multiple access via Approximate Message
# Passing" (the paper that now has an arXiv link), so find that specific
# paragraph first and restrict the Find/Replace to its Range.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Approximate Message Passing*") {
        $r = $p.Range
        $r.Find.Execute(", draft available on request", $true, $false, $false,
                         $false, $false, $true, 1, $false, "", 2)
    }
}
